# Update transition-probability matrix values in the single worksheet
# to reflect the recomputed "team specific time" statistics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1581196581196581
$ws.Range("C2").Value = 0.6324786324786325
$ws.Range("J2").Value = 0.01282051282051282
$ws.Range("P2").Value = 0.1324786324786325
$ws.Range("S2").Value = 0.0641025641025641
$ws.Range("B3").Value = 0.01298701298701299
$ws.Range("C3").Value = 0.03246753246753246
$ws.Range("J3").Value = 0.01948051948051948
$ws.Range("P3").Value = 0.7987012987012987
$ws.Range("S3").Value = 0.1363636363636364
$ws.Range("J4").Value = 0.03333333333333333
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.3
$ws.Range("B6").Value = 0.0735930735930736
$ws.Range("D6").Value = 0.01298701298701299
$ws.Range("E6").Value = 0.004329004329004329
$ws.Range("F6").Value = 0.02597402597402598
$ws.Range("J6").Value = 0.2987012987012987
$ws.Range("O6").Value = 0.01298701298701299
$ws.Range("Q6").Value = 0.1558441558441558
$ws.Range("R6").Value = 0.1038961038961039
$ws.Range("S6").Value = 0.3116883116883117
$ws.Range("B7").Value = 0.08994708994708994
$ws.Range("D7").Value = 0.01587301587301587
$ws.Range("F7").Value = 0.06349206349206349
$ws.Range("J7").Value = 0.1164021164021164
$ws.Range("O7").Value = 0.03174603174603174
$ws.Range("Q7").Value = 0.201058201058201
$ws.Range("R7").Value = 0.07407407407407407
$ws.Range("S7").Value = 0.4074074074074074
$ws.Range("B8").Value = 0.07959183673469387
$ws.Range("D8").Value = 0.01224489795918367
$ws.Range("F8").Value = 0.06122448979591837
$ws.Range("J8").Value = 0.1081632653061225
$ws.Range("O8").Value = 0.03061224489795918
$ws.Range("Q8").Value = 0.1897959183673469
$ws.Range("R8").Value = 0.1040816326530612
$ws.Range("S8").Value = 0.4142857142857143
$ws.Range("B9").Value = 0.05291005291005291
$ws.Range("F9").Value = 0.0582010582010582
$ws.Range("J9").Value = 0.1216931216931217
$ws.Range("O9").Value = 0.01058201058201058
$ws.Range("Q9").Value = 0.1746031746031746
$ws.Range("R9").Value = 0.1164021164021164
$ws.Range("S9").Value = 0.4656084656084656
$ws.Range("B10").Value = 0.07607192254495158
$ws.Range("D10").Value = 0.01313969571230982
$ws.Range("E10").Value = 0.0006915629322268327
$ws.Range("F10").Value = 0.0698478561549101
$ws.Range("J10").Value = 0.1044260027662517
$ws.Range("O10").Value = 0.01383125864453665
$ws.Range("Q10").Value = 0.2593360995850623
$ws.Range("R10").Value = 0.08506224066390042
$ws.Range("S10").Value = 0.3775933609958506
$ws.Range("G11").Value = 0.1388888888888889
$ws.Range("J11").Value = 0.0798611111111111
$ws.Range("K11").Value = 0.2083333333333333
$ws.Range("L11").Value = 0.5659722222222222
$ws.Range("S11").Value = 0.006944444444444444
$ws.Range("G12").Value = 0.7455621301775148
$ws.Range("J12").Value = 0.2189349112426036
$ws.Range("K12").Value = 0.005917159763313609
$ws.Range("L12").Value = 0.02958579881656805
$ws.Range("G13").Value = 0.5740740740740741
$ws.Range("J13").Value = 0.4259259259259259
$ws.Range("F15").Value = 0.01282051282051282
$ws.Range("H15").Value = 0.1623931623931624
$ws.Range("I15").Value = 0.05128205128205128
$ws.Range("J15").Value = 0.405982905982906
$ws.Range("K15").Value = 0.04273504273504274
$ws.Range("M15").Value = 0.01282051282051282
$ws.Range("O15").Value = 0.04700854700854701
$ws.Range("S15").Value = 0.264957264957265
$ws.Range("F16").Value = 0.02352941176470588
$ws.Range("H16").Value = 0.2176470588235294
$ws.Range("I16").Value = 0.07058823529411765
$ws.Range("J16").Value = 0.4588235294117647
$ws.Range("K16").Value = 0.09411764705882353
$ws.Range("M16").Value = 0.01764705882352941
$ws.Range("O16").Value = 0.03529411764705882
$ws.Range("S16").Value = 0.08235294117647059
$ws.Range("F17").Value = 0.02112676056338028
$ws.Range("H17").Value = 0.1901408450704225
$ws.Range("I17").Value = 0.07570422535211267
$ws.Range("J17").Value = 0.4471830985915493
$ws.Range("K17").Value = 0.09330985915492958
$ws.Range("M17").Value = 0.02816901408450704
$ws.Range("N17").Value = 0.00352112676056338
$ws.Range("O17").Value = 0.06514084507042253
$ws.Range("S17").Value = 0.07570422535211267
$ws.Range("F18").Value = 0.03017241379310345
$ws.Range("H18").Value = 0.2025862068965517
$ws.Range("I18").Value = 0.08620689655172414
$ws.Range("J18").Value = 0.4439655172413793
$ws.Range("K18").Value = 0.08620689655172414
$ws.Range("M18").Value = 0.01724137931034483
$ws.Range("O18").Value = 0.0603448275862069
$ws.Range("S18").Value = 0.07327586206896551
$ws.Range("F19").Value = 0.0136327185244587
$ws.Range("H19").Value = 0.210906174819567
$ws.Range("I19").Value = 0.08259823576583801
$ws.Range("J19").Value = 0.4186046511627907
$ws.Range("K19").Value = 0.1018444266238974
$ws.Range("M19").Value = 0.02405773857257418
$ws.Range("N19").Value = 0.0008019246190858059
$ws.Range("O19").Value = 0.07056936647955092
$ws.Range("S19").Value = 0.07698476343223737
